# refactor currency conversion, now explicit source and target amounts
#
# The "currency_conversions" sheet used to track a single "foreign_amount"
# (the amount in the source currency) together with "source_fees" only.
# This splits the conversion into an explicit source side (source_amount,
# source_fees, source_currency) and target side (target_amount, target_fees,
# target_currency), followed by the free-text comment column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("currency_conversions")

# Old layout: date | foreign_amount | source_fees | source_currency | target_currency | comment
# New layout: date | source_amount  | source_fees | source_currency | target_amount | target_fees | target_currency | comment
# Insert two fresh columns right before the old "target_currency" column (E)
# to hold the new target_amount / target_fees fields.
$ws.Range("E1:F1").EntireColumn.Insert()

# "foreign_amount" becomes "source_amount" (same column/values, renamed header).
$ws.Range("B1").Value = "source_amount"

# Populate the two newly inserted columns.
$ws.Range("E1").Value = "target_amount"
$ws.Range("F1").Value = "target_fees"

$ws.Range("E2").Value = -1
$ws.Range("F2").Value = 0

$ws.Range("E3").Value = -1
$ws.Range("F3").Value = 0

# currency_conversions becomes the active sheet/tab.
$ws.Activate()
